$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9459
$ws.Range("E2").Value = 8367
$ws.Range("F2").Value = 0.8845543926419284
$ws.Range("G2").Value = 0.8827811774636
$ws.Range("H2").Value = 0.09590404283678339
$ws.Range("I2").Value = 0.08466228385897516
$ws.Range("J2").Value = 40698145.31691629
$ws.Range("K2").Value = 14154240.60511414
$ws.Range("M2").Value = 14154240.60511414
$ws.Range("N2").Value = 54852385.92203042
$ws.Range("O2").Value = 801445750.0172
$ws.Range("P2").Value = 783745943.0132
$ws.Range("Q2").Value = 0.01766088422679935
$ws.Range("R2").Value = 0.01805973061971659

# Row 3
$ws.Range("D3").Value = 9637
$ws.Range("E3").Value = 8559
$ws.Range("F3").Value = 0.8881394624883262
$ws.Range("G3").Value = 0.8854748603351955
$ws.Range("H3").Value = 0.09414955439058634
$ws.Range("I3").Value = 0.08336706352462532
$ws.Range("J3").Value = 42370594.43989093
$ws.Range("K3").Value = 14720548.15968467
$ws.Range("M3").Value = 14720548.15968467
$ws.Range("N3").Value = 57091142.59957561
$ws.Range("O3").Value = 836286487.948328
$ws.Range("P3").Value = 818806311.924258
$ws.Range("Q3").Value = 0.01760227908954833
$ws.Range("R3").Value = 0.0179780589686592

# Row 4
$ws.Range("D4").Value = 9838
$ws.Range("E4").Value = 8733
$ws.Range("F4").Value = 0.8876804228501728
$ws.Range("G4").Value = 0.8858794887401096
$ws.Range("H4").Value = 0.09295468746118477
$ws.Range("I4").Value = 0.08234665100411104
$ws.Range("J4").Value = 44324690.6429288
$ws.Range("K4").Value = 15377674.13053978
$ws.Range("M4").Value = 15377674.13053978
$ws.Range("N4").Value = 59702364.77346857
$ws.Range("O4").Value = 875462780.5151056
$ws.Range("P4").Value = 858013832.5091518
$ws.Range("Q4").Value = 0.01756519462939572
$ws.Range("R4").Value = 0.01792240818026178

# Row 5
$ws.Range("D5").Value = 10033
$ws.Range("E5").Value = 8917
$ws.Range("F5").Value = 0.8887670686733778
$ws.Range("G5").Value = 0.8869106823154963
$ws.Range("H5").Value = 0.0915457106288796
$ws.Range("I5").Value = 0.08119286867691659
$ws.Range("J5").Value = 46252928.80369589
$ws.Range("K5").Value = 16008719.50273435
$ws.Range("M5").Value = 16008719.50273435
$ws.Range("N5").Value = 62261648.30643024
$ws.Range("O5").Value = 914819837.8515847
$ws.Range("P5").Value = 897333732.3886114
$ws.Range("Q5").Value = 0.01749931389805685
$ws.Range("R5").Value = 0.01784031840653172

# Row 6
$ws.Range("D6").Value = 10220
$ws.Range("E6").Value = 9059
$ws.Range("F6").Value = 0.886399217221135
$ws.Range("G6").Value = 0.8834601131265848
$ws.Range("H6").Value = 0.09057818461513951
$ws.Range("I6").Value = 0.08002221322689183
$ws.Range("J6").Value = 48283872.6099349
$ws.Range("K6").Value = 16653146.74602686
$ws.Range("M6").Value = 16653146.74602686
$ws.Range("N6").Value = 64937019.35596175
$ws.Range("O6").Value = 955659491.8934135
$ws.Range("P6").Value = 938067666.020232
$ws.Range("Q6").Value = 0.01742581629470616
$ws.Range("R6").Value = 0.01775260714046154

# Row 7
$ws.Range("D7").Value = 9444
$ws.Range("E7").Value = 8346
$ws.Range("F7").Value = 0.8837357052096569
$ws.Range("G7").Value = 0.8805655201519308
$ws.Range("H7").Value = 0.0970076441535767
$ws.Range("I7").Value = 0.08542158663280766
$ws.Range("J7").Value = 41076824.4099903
$ws.Range("K7").Value = 14343580.15165115
$ws.Range("M7").Value = 14343580.15165115
$ws.Range("N7").Value = 55420404.56164145
$ws.Range("O7").Value = 798913812.9572001
$ws.Range("P7").Value = 781214005.9532001
$ws.Range("Q7").Value = 0.01795385174097568
$ws.Range("R7").Value = 0.01836062851196554

# Row 8
$ws.Range("D8").Value = 9640
$ws.Range("E8").Value = 8568
$ws.Range("F8").Value = 0.8887966804979253
$ws.Range("G8").Value = 0.8864059590316573
$ws.Range("H8").Value = 0.09542172038036506
$ws.Range("I8").Value = 0.08458238156620813
$ws.Range("J8").Value = 43142786.81432747
$ws.Range("K8").Value = 15106644.34690293
$ws.Range("M8").Value = 15106644.34690293
$ws.Range("N8").Value = 58249431.1612304
$ws.Range("O8").Value = 837045574.001528
$ws.Range("P8").Value = 819565397.977458
$ws.Range("Q8").Value = 0.01804757687766635
$ws.Range("R8").Value = 0.01843250628221183

# Row 9
$ws.Range("D9").Value = 9845
$ws.Range("E9").Value = 8756
$ws.Range("F9").Value = 0.8893854748603351
$ws.Range("G9").Value = 0.888212619192534
$ws.Range("H9").Value = 0.09412919338671218
$ws.Range("I9").Value = 0.08360673740049218
$ws.Range("J9").Value = 45118047.21266638
$ws.Range("K9").Value = 15774352.41540857
$ws.Range("M9").Value = 15774352.41540857
$ws.Range("N9").Value = 60892399.62807495
$ws.Range("O9").Value = 874986044.7346259
$ws.Range("P9").Value = 857537096.7286721
$ws.Range("Q9").Value = 0.01802811886010453
$ws.Range("R9").Value = 0.0183949504640493

# Row 10
$ws.Range("D10").Value = 10031
$ws.Range("E10").Value = 8913
$ws.Range("F10").Value = 0.8885455089223407
$ws.Range("G10").Value = 0.8865128307141437
$ws.Range("H10").Value = 0.09321308070362111
$ws.Range("I10").Value = 0.08263459203415308
$ws.Range("J10").Value = 47213184.81612386
$ws.Range("K10").Value = 16488847.50894834
$ws.Range("M10").Value = 16488847.50894834
$ws.Range("N10").Value = 63702032.3250722
$ws.Range("O10").Value = 914432692.5229203
$ws.Range("P10").Value = 896946587.059947
$ws.Range("Q10").Value = 0.01803177822028169
$ws.Range("R10").Value = 0.01838331038528866

# Row 11
$ws.Range("D11").Value = 10225
$ws.Range("E11").Value = 9080
$ws.Range("F11").Value = 0.8880195599022005
$ws.Range("G11").Value = 0.8855080944021845
$ws.Range("H11").Value = 0.09213732329015105
$ws.Range("I11").Value = 0.08158834556997968
$ws.Range("J11").Value = 49406731.03858929
$ws.Range("K11").Value = 17214575.96035406
$ws.Range("M11").Value = 17214575.96035406
$ws.Range("N11").Value = 66621306.99894334
$ws.Range("O11").Value = 956127888.2484893
$ws.Range("P11").Value = 938536062.3753077
$ws.Range("Q11").Value = 0.01800447008390172
$ws.Range("R11").Value = 0.01834194406636469

# Row 12
$ws.Range("D12").Value = 9455
$ws.Range("E12").Value = 8375
$ws.Range("F12").Value = 0.8857747223691169
$ws.Range("G12").Value = 0.8836252373918548
$ws.Range("H12").Value = 0.09689384242872673
$ws.Range("I12").Value = 0.08561784451789264
$ws.Range("J12").Value = 41234906.25223832
$ws.Range("K12").Value = 14421586.78406116
$ws.Range("M12").Value = 14421586.78406116
$ws.Range("N12").Value = 55656493.03629947
$ws.Range("O12").Value = 800568137.6472001
$ws.Range("P12").Value = 782868330.6432
$ws.Range("Q12").Value = 0.018014190305456
$ws.Range("R12").Value = 0.01842147168248901

# Row 13
$ws.Range("D13").Value = 9647
$ws.Range("E13").Value = 8557
$ws.Range("F13").Value = 0.8870115061677205
$ws.Range("G13").Value = 0.8852679495137595
$ws.Range("H13").Value = 0.1020477647857244
$ws.Range("I13").Value = 0.09033961548432071
$ws.Range("J13").Value = 48109582.45102569
$ws.Range("K13").Value = 17588976.84787663
$ws.Range("M13").Value = 17588976.84787663
$ws.Range("N13").Value = 65698559.29890232
$ws.Range("O13").Value = 838121313.564728
$ws.Range("P13").Value = 820641137.540658
$ws.Range("Q13").Value = 0.02098619443653873
$ws.Range("R13").Value = 0.0214332136707017

# Row 14
$ws.Range("D14").Value = 9840
$ws.Range("E14").Value = 8723
$ws.Range("F14").Value = 0.8864837398373984
$ws.Range("G14").Value = 0.8848650841955772
$ws.Range("H14").Value = 0.1062957128902854
$ws.Range("I14").Value = 0.09405736493629131
$ws.Range("J14").Value = 54616763.50503325
$ws.Range("K14").Value = 20522613.28469532
$ws.Range("M14").Value = 20522613.28469532
$ws.Range("N14").Value = 75139376.78972858
$ws.Range("O14").Value = 876530452.3962009
$ws.Range("P14").Value = 859081504.3902471
$ws.Range("Q14").Value = 0.02341346296479713
$ws.Range("R14").Value = 0.02388901772394893

# Row 15
$ws.Range("D15").Value = 10032
$ws.Range("E15").Value = 8909
$ws.Range("F15").Value = 0.8880582137161085
$ws.Range("G15").Value = 0.8861149791127909
$ws.Range("H15").Value = 0.1093961613107875
$ws.Range("I15").Value = 0.09693757719492795
$ws.Range("J15").Value = 60672123.38640694
$ws.Range("K15").Value = 23217186.5988863
$ws.Range("M15").Value = 23217186.5988863
$ws.Range("N15").Value = 83889309.98529324
$ws.Range("O15").Value = 914302919.6242424
$ws.Range("P15").Value = 896816814.1612692
$ws.Range("Q15").Value = 0.02539331998242774
$ws.Range("R15").Value = 0.02588843812055389

# Row 16
$ws.Range("D16").Value = 10235
$ws.Range("E16").Value = 9057
$ws.Range("F16").Value = 0.884904738641915
$ws.Range("G16").Value = 0.8832650672908133
$ws.Range("H16").Value = 0.1090150774465601
$ws.Range("I16").Value = 0.0962892097165491
$ws.Range("J16").Value = 63804904.40937157
$ws.Range("K16").Value = 24412498.54468551
$ws.Range("M16").Value = 24412498.54468551
$ws.Range("N16").Value = 88217402.95405708
$ws.Range("O16").Value = 955199529.3546511
$ws.Range("P16").Value = 937607703.4814696
$ws.Range("Q16").Value = 0.02555748594346461
$ws.Range("R16").Value = 0.02603700721958497
